# Auto-generated edit script to apply profit-table updates across multiple sheets
# (values sourced from the commit's scheduled-runner price refresh)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(58, 8).Value = 202.83333
$ws.Cells.Item(58, 9).Value = 202.83333
$ws.Cells.Item(58, 11).Value = 608.49999
$ws.Cells.Item(58, 13).Value = -458.49999

$ws.Cells.Item(98, 8).Value = 1597.1428
$ws.Cells.Item(98, 9).Value = 1608.3077
$ws.Cells.Item(98, 10).Value = 1452
$ws.Cells.Item(98, 11).Value = 1608.3077
$ws.Cells.Item(98, 12).Value = 1452
$ws.Cells.Item(98, 13).Value = -110.3077000000001
$ws.Cells.Item(98, 14).Value = -4448

$ws.Cells.Item(103, 8).Value = 1230.25
$ws.Cells.Item(103, 9).Value = 973.6667
$ws.Cells.Item(103, 10).Value = 2000
$ws.Cells.Item(103, 11).Value = 2921.0001
$ws.Cells.Item(103, 12).Value = 6000
$ws.Cells.Item(103, 13).Value = -2335.0001
$ws.Cells.Item(103, 14).Value = -7172

$ws.Cells.Item(112, 8).Value = 5250.5293
$ws.Cells.Item(112, 10).Value = 5380.8486
$ws.Cells.Item(112, 12).Value = 16142.5458
$ws.Cells.Item(112, 14).Value = -18358.5458

$ws.Cells.Item(122, 8).Value = 1597.1428
$ws.Cells.Item(122, 9).Value = 1608.3077
$ws.Cells.Item(122, 10).Value = 1452
$ws.Cells.Item(122, 11).Value = 4824.9231
$ws.Cells.Item(122, 12).Value = 4356
$ws.Cells.Item(122, 13).Value = -2374.9231
$ws.Cells.Item(122, 14).Value = -9256

$ws.Cells.Item(132, 8).Value = 67005.94500000001
$ws.Cells.Item(132, 9).Value = 72327.39
$ws.Cells.Item(132, 11).Value = 216982.17
$ws.Cells.Item(132, 13).Value = -214452.17

$ws.Cells.Item(135, 8).Value = 1702.0303
$ws.Cells.Item(135, 10).Value = 3003.923
$ws.Cells.Item(135, 12).Value = 27035.307
$ws.Cells.Item(135, 14).Value = -32105.307

$ws.Cells.Item(137, 8).Value = 1626324.1
$ws.Cells.Item(137, 10).Value = 3610860.5
$ws.Cells.Item(137, 12).Value = 10832581.5
$ws.Cells.Item(137, 14).Value = -10837681.5

$ws.Cells.Item(140, 8).Value = 40590
$ws.Cells.Item(140, 10).Value = 40590
$ws.Cells.Item(140, 12).Value = 40590
$ws.Cells.Item(140, 14).Value = -50950

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 22991.72
$ws.Cells.Item(32, 9).Value = 23452.98
$ws.Cells.Item(32, 11).Value = 23452.98
$ws.Cells.Item(32, 13).Value = -23165.98

$ws.Cells.Item(61, 8).Value = 6667705
$ws.Cells.Item(61, 9).Value = 6667705
$ws.Cells.Item(61, 11).Value = 6667705
$ws.Cells.Item(61, 13).Value = -6667493

$ws.Cells.Item(63, 8).Value = 13656.889
$ws.Cells.Item(63, 9).Value = 1899.8334
$ws.Cells.Item(63, 10).Value = 37171
$ws.Cells.Item(63, 11).Value = 1899.8334
$ws.Cells.Item(63, 12).Value = 37171
$ws.Cells.Item(63, 13).Value = -1213.8334
$ws.Cells.Item(63, 14).Value = -38543

$ws.Cells.Item(66, 8).Value = 13656.889
$ws.Cells.Item(66, 9).Value = 1899.8334
$ws.Cells.Item(66, 10).Value = 37171
$ws.Cells.Item(66, 11).Value = 9499.166999999999
$ws.Cells.Item(66, 12).Value = 185855
$ws.Cells.Item(66, 13).Value = -6067.166999999999
$ws.Cells.Item(66, 14).Value = -192719

$ws.Cells.Item(122, 8).Value = 740.7273
$ws.Cells.Item(122, 9).Value = 740.7273
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 2222.1819
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = 227.8181
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(132, 8).Value = 386670.28
$ws.Cells.Item(132, 9).Value = 436453.66
$ws.Cells.Item(132, 11).Value = 1309360.98
$ws.Cells.Item(132, 13).Value = -1306830.98

$ws.Cells.Item(133, 8).Value = 97999.5
$ws.Cells.Item(133, 10).Value = 97999.5
$ws.Cells.Item(133, 12).Value = 97999.5
$ws.Cells.Item(133, 14).Value = -103059.5

$ws.Cells.Item(136, 8).Value = 6667705
$ws.Cells.Item(136, 9).Value = 6667705
$ws.Cells.Item(136, 11).Value = 20003115
$ws.Cells.Item(136, 13).Value = -20000565

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1760.45
$ws.Cells.Item(86, 9).Value = 1760.45
$ws.Cells.Item(86, 11).Value = 1760.45
$ws.Cells.Item(86, 13).Value = -637.45

$ws.Cells.Item(89, 8).Value = 1760.45
$ws.Cells.Item(89, 9).Value = 1760.45
$ws.Cells.Item(89, 11).Value = 8802.25
$ws.Cells.Item(89, 13).Value = -3186.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 18836.385
$ws.Cells.Item(31, 9).Value = 8259.385
$ws.Cells.Item(31, 10).Value = 29413.385
$ws.Cells.Item(31, 11).Value = 8259.385
$ws.Cells.Item(31, 12).Value = 29413.385
$ws.Cells.Item(31, 13).Value = -7964.385
$ws.Cells.Item(31, 14).Value = -30003.385

$ws.Cells.Item(34, 8).Value = 18836.385
$ws.Cells.Item(34, 9).Value = 8259.385
$ws.Cells.Item(34, 10).Value = 29413.385
$ws.Cells.Item(34, 11).Value = 8259.385
$ws.Cells.Item(34, 12).Value = 29413.385
$ws.Cells.Item(34, 13).Value = -8057.385
$ws.Cells.Item(34, 14).Value = -29817.385

$ws.Cells.Item(107, 8).Value = 1214.125
$ws.Cells.Item(107, 9).Value = 1054.8889
$ws.Cells.Item(107, 11).Value = 1054.8889
$ws.Cells.Item(107, 13).Value = 865.1111000000001

$ws.Cells.Item(134, 8).Value = 5537.543
$ws.Cells.Item(134, 9).Value = 5611.853
$ws.Cells.Item(134, 11).Value = 16835.559
$ws.Cells.Item(134, 13).Value = -14300.559

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 666.6667
$ws.Cells.Item(5, 9).Value = 500
$ws.Cells.Item(5, 10).Value = 1000
$ws.Cells.Item(5, 11).Value = 1500
$ws.Cells.Item(5, 12).Value = 3000
$ws.Cells.Item(5, 13).Value = -1388
$ws.Cells.Item(5, 14).Value = -3224

$ws.Cells.Item(46, 8).Value = 6000
$ws.Cells.Item(46, 10).Value = 6000
$ws.Cells.Item(46, 12).Value = 18000
$ws.Cells.Item(46, 14).Value = -18182

$ws.Cells.Item(58, 8).Value = 663
$ws.Cells.Item(58, 9).Value = 670
$ws.Cells.Item(58, 10).Value = 649
$ws.Cells.Item(58, 11).Value = 2010
$ws.Cells.Item(58, 12).Value = 1947
$ws.Cells.Item(58, 13).Value = -1882
$ws.Cells.Item(58, 14).Value = -2203

$ws.Cells.Item(109, 8).Value = 1931.25
$ws.Cells.Item(109, 9).Value = 1408.3334
$ws.Cells.Item(109, 10).Value = 3500
$ws.Cells.Item(109, 11).Value = 4225.0002
$ws.Cells.Item(109, 12).Value = 10500
$ws.Cells.Item(109, 13).Value = -3185.0002
$ws.Cells.Item(109, 14).Value = -12580

$ws.Cells.Item(132, 8).Value = 1000
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1000
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).Value = 9000
$ws.Cells.Item(132, 13).ClearContents()
$ws.Cells.Item(132, 14).Value = -14060

$ws.Cells.Item(135, 8).Value = 666.6667
$ws.Cells.Item(135, 9).Value = 500
$ws.Cells.Item(135, 10).Value = 1000
$ws.Cells.Item(135, 11).Value = 4500
$ws.Cells.Item(135, 12).Value = 9000
$ws.Cells.Item(135, 13).Value = -1965
$ws.Cells.Item(135, 14).Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 4731.3335
$ws.Cells.Item(70, 9).Value = 4747
$ws.Cells.Item(70, 11).Value = 4747
$ws.Cells.Item(70, 13).Value = -4477

$ws.Cells.Item(73, 8).Value = 4731.3335
$ws.Cells.Item(73, 9).Value = 4747
$ws.Cells.Item(73, 11).Value = 4747
$ws.Cells.Item(73, 13).Value = -3811

$ws.Cells.Item(74, 8).Value = 81282.75

$ws.Cells.Item(77, 8).Value = 81282.75

$ws.Cells.Item(122, 8).Value = 46963.08
$ws.Cells.Item(122, 10).Value = 9740
$ws.Cells.Item(122, 12).Value = 29220
$ws.Cells.Item(122, 14).Value = -34120

$ws.Cells.Item(126, 8).Value = 1391722.1
$ws.Cells.Item(126, 9).Value = 2383334
$ws.Cells.Item(126, 10).Value = 3465.6
$ws.Cells.Item(126, 11).Value = 7150002
$ws.Cells.Item(126, 12).Value = 10396.8
$ws.Cells.Item(126, 13).Value = -7147532
$ws.Cells.Item(126, 14).Value = -15336.8

$ws.Cells.Item(132, 8).Value = 22005232
$ws.Cells.Item(132, 9).Value = 28115946
$ws.Cells.Item(132, 10).Value = 6661.1
$ws.Cells.Item(132, 11).Value = 84347838
$ws.Cells.Item(132, 12).Value = 19983.3
$ws.Cells.Item(132, 13).Value = -84345308
$ws.Cells.Item(132, 14).Value = -25043.3

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 3865.3333
$ws.Cells.Item(7, 9).Value = 3598.4119
$ws.Cells.Item(7, 11).Value = 3598.4119
$ws.Cells.Item(7, 13).Value = -3486.4119

$ws.Cells.Item(93, 8).Value = 2999
$ws.Cells.Item(93, 9).Value = 2999
$ws.Cells.Item(93, 11).Value = 2999
$ws.Cells.Item(93, 13).Value = -1751

$ws.Cells.Item(100, 8).Value = 9585.857
$ws.Cells.Item(100, 9).Value = 2199.7778
$ws.Cells.Item(100, 11).Value = 2199.7778
$ws.Cells.Item(100, 13).Value = -1658.7778

$ws.Cells.Item(126, 8).Value = 3865.3333
$ws.Cells.Item(126, 9).Value = 3598.4119
$ws.Cells.Item(126, 11).Value = 10795.2357
$ws.Cells.Item(126, 13).Value = -8325.235700000001

$ws.Cells.Item(132, 8).Value = 2680767.2
$ws.Cells.Item(132, 9).Value = 3869854
$ws.Cells.Item(132, 10).Value = 5321.75
$ws.Cells.Item(132, 11).Value = 11609562
$ws.Cells.Item(132, 12).Value = 15965.25
$ws.Cells.Item(132, 13).Value = -11607032
$ws.Cells.Item(132, 14).Value = -21025.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(101, 8).Value = 31750
$ws.Cells.Item(101, 10).Value = 31750
$ws.Cells.Item(101, 12).Value = 31750
$ws.Cells.Item(101, 14).Value = -38240

$ws.Cells.Item(104, 8).Value = 50185
$ws.Cells.Item(104, 10).Value = 50185
$ws.Cells.Item(104, 12).Value = 50185
$ws.Cells.Item(104, 14).Value = -57173

$ws.Cells.Item(107, 8).Value = 1371.3334
$ws.Cells.Item(107, 9).Value = 746.0454999999999
$ws.Cells.Item(107, 10).Value = 2621.9092
$ws.Cells.Item(107, 11).Value = 2238.1365
$ws.Cells.Item(107, 12).Value = 7865.7276
$ws.Cells.Item(107, 13).Value = -318.1364999999996
$ws.Cells.Item(107, 14).Value = -11705.7276

$ws.Cells.Item(113, 8).Value = 2814.238
$ws.Cells.Item(113, 10).Value = 4042.2307
$ws.Cells.Item(113, 12).Value = 12126.6921
$ws.Cells.Item(113, 14).Value = -16466.6921
